$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 404
$wsExpo.Range("F4").Value = 5082
$wsExpo.Range("F5").Value = 43
$wsExpo.Range("F6").Value = 41
$wsExpo.Range("F8").Value = 501

# Sheet "全部类型" (all types) - same events, same updated values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 404
$wsAll.Range("F4").Value = 5082
$wsAll.Range("F6").Value = 43
$wsAll.Range("F7").Value = 41
$wsAll.Range("F10").Value = 501
